$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203425049781799
$ws.Range("B1").Value = 1.394992709159851
$ws.Range("C1").Value = 3.611523389816284
$ws.Range("D1").Value = 3.645560503005981
$ws.Range("E1").Value = 0.9987479448318481
